# Updating barcode_offset and umi_offset
#
# 1. barcode_offset sheet: the "0" offset option becomes "0,38,76", and a new
#    option "10,48,86" is appended.
# 2. umi_offset sheet: a new option "1" is inserted before the existing
#    options ("Not applicable", "16").
# 3. The ATACseq sheet's data-validation ranges that reference those two
#    lookup sheets are widened to include the extra row.
# 4. The .metadata sheet's pav:createdOn timestamp is refreshed.

$wb = $excel.ActiveWorkbook

# --- barcode_offset --------------------------------------------------
$boSheet = $wb.Worksheets.Item("barcode_offset")
$boSheet.Range("A1").Value = "0,38,76"
$boSheet.Range("A5").Value = "10,48,86"

# --- umi_offset --------------------------------------------------------
$uoSheet = $wb.Worksheets.Item("umi_offset")
$uoSheet.Rows.Item(1).Insert()
$uoSheet.Range("A1").NumberFormat = "@"
$uoSheet.Range("A1").Value = "1"

# --- ATACseq data validation ranges -------------------------------------
$mainSheet = $wb.Worksheets.Item("ATACseq")
$mainSheet.Range("O2:O1001").Validation.Modify(3, 1, 1, "'barcode_offset'!`$A`$1:`$A`$5")
$mainSheet.Range("R2:R1001").Validation.Modify(3, 1, 1, "'umi_offset'!`$A`$1:`$A`$3")

# --- .metadata: pav:createdOn timestamp ---------------------------------
$metaSheet = $wb.Worksheets.Item(".metadata")
$metaSheet.Range("C2").Value = "2023-10-31T13:53:10-07:00"
